$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Column J ("EndingPara" / "list,int" / "结局参数补正") is being dropped,
# and the old column M ("UseCondition" / "EventCondition" / "使用条件")
# takes its place at J (with its row-2 type text updated), then the
# now-empty column M is removed so K ("ResourceTags") and L
# ("ResourceImage") stay put.
# ------------------------------------------------------------------
$ws.Range("J1").Value2 = "UseCondition"
$ws.Range("J2").Value2 = "(list#sep=;), EventCondition#sep=,"
$ws.Range("J3").Value2 = "使用条件"

$ws.Columns("M").Delete() | Out-Null

# Highlight the updated type cell (J2) with a solid yellow fill.
$ws.Range("J2").Interior.Color = 65535

# Column width tweaks for G, I, J (new widths after the column shuffle).
$ws.Columns("G").ColumnWidth = 13.43
$ws.Columns("I").ColumnWidth = 13.71
$ws.Columns("J").ColumnWidth = 11.29

# Move the active selection (cosmetic, mirrors the saved view state).
$ws.Range("J9").Select() | Out-Null
